# Add a "2022-Q1" sheet (repurposing the old "总计" sheet's slot) and
# append a fresh "总计" sheet after it, carrying forward the original
# summary table plus the new 2022-Q1 row.

$wb = $excel.ActiveWorkbook

$sheetQ4   = $wb.Worksheets.Item(4)   # "2021-Q4" - same 9-row fund layout as the new sheet needs (minus one row)
$sheetSlot = $wb.Worksheets.Item(5)   # currently "总计" - becomes "2022-Q1"

# --- 1. Create the new (placeholder-named) sheet right after the current
#        slot, and copy the ORIGINAL "总计" formatting onto it before that
#        sheet's own content gets overwritten below. The rename to "总计"
#        itself has to wait until the old slot has been renamed away from
#        "总计" (Excel won't allow two sheets with the same name).
$sheetTotal = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheetSlot)

# (copy the header-row and index-column formatting separately so A1 - which
#  never holds a value on these sheets - never gets materialised as an
#  empty cell element)
$sheetSlot.Range("B1:D1").Copy()
$sheetTotal.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$sheetSlot.Range("A2:A5").Copy()
$sheetTotal.Range("A2").PasteSpecial(-4122)
$sheetSlot.Range("A5").Copy()
$sheetTotal.Range("A6").PasteSpecial(-4122)   # extend the row-5 look to the new row 6

# --- 2. Re-purpose the slot sheet: rename + restyle to match the other
#        fund-holding sheets (9 data rows, 8 columns A:H), then fill values.
$sheetSlot.Name = "2022-Q1"
$sheetTotal.Name = "总计"

$sheetQ4.Range("B1:H1").Copy()
$sheetSlot.Range("B1").PasteSpecial(-4122)    # xlPasteFormats
$sheetQ4.Range("A2:A8").Copy()
$sheetSlot.Range("A2").PasteSpecial(-4122)
$sheetQ4.Range("A8").Copy()
$sheetSlot.Range("A9").PasteSpecial(-4122)    # extend row-8 styling down to the extra row 9

function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Header row
$sheetSlot.Range("B1").Value = "基金代码"
$sheetSlot.Range("C1").Value = "基金名称"
$sheetSlot.Range("D1").Value = "基金规模"
$sheetSlot.Range("E1").Value = "股票总仓位"
$sheetSlot.Range("F1").Value = "仓位占比"
$sheetSlot.Range("G1").Value = "持有市值(亿元)"
$sheetSlot.Range("H1").Value = "仓位排名"

# Data rows (fund code / fund name / regulated-as-text numeric columns / rank as real number)
$rows = @(
    @{A=0; B="100056"; C="富国低碳环保混合";                 D="27.82"; E="82.05"; F="2.69"; G="0.7484"; H=9},
    @{A=1; B="001476"; C="中银智能制造股票";                 D="19.50"; E="82.72"; F="3.49"; G="0.6806"; H=3},
    @{A=2; B="163807"; C="中银优选混合";                     D="30.26"; E="72.71"; F="2.06"; G="0.6234"; H=5},
    @{A=3; B="009693"; C="富国积极成长一年定期开放混合";      D="17.82"; E="98.74"; F="2.67"; G="0.4758"; H=9},
    @{A=4; B="011212"; C="富国稳健策略6个月持有期混合A";      D="12.61"; E="83.04"; F="3.00"; G="0.3783"; H=8},
    @{A=5; B="009379"; C="中银成长优选股票";                 D="2.18";  E="81.34"; F="5.26"; G="0.1147"; H=1},
    @{A=6; B="011213"; C="富国稳健策略6个月持有期混合C";      D="1.14";  E="83.04"; F="3.00"; G="0.0342"; H=8},
    @{A=7; B="320016"; C="诺安多策略混合";                   D="0.19";  E="80.02"; F="3.85"; G="0.0073"; H=6}
)

$r = 2
foreach ($row in $rows) {
    $sheetSlot.Cells.Item($r, 1).Value = $row.A
    Set-TextCell $sheetSlot $r 2 $row.B
    Set-TextCell $sheetSlot $r 3 $row.C
    Set-TextCell $sheetSlot $r 4 $row.D
    Set-TextCell $sheetSlot $r 5 $row.E
    Set-TextCell $sheetSlot $r 6 $row.F
    Set-TextCell $sheetSlot $r 7 $row.G
    $sheetSlot.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# --- 3. Fill the new "总计" sheet: 2022-Q1 row on top, followed by the
#        pre-existing quarterly summary rows (shifted down by one).
$sheetTotal.Range("B1").Value = "日期"
$sheetTotal.Range("C1").Value = "持有数量(只)"
$sheetTotal.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @{A=0; B="2022-Q1"; C=8;  D=3.06},
    @{A=1; B="2021-Q4"; C=7;  D=3.6},
    @{A=2; B="2021-Q3"; C=17; D=3.9},
    @{A=3; B="2021-Q2"; C=25; D=7.15},
    @{A=4; B="2021-Q1"; C=7;  D=0.53}
)

$r = 2
foreach ($row in $totalRows) {
    $sheetTotal.Cells.Item($r, 1).Value = $row.A
    $sheetTotal.Cells.Item($r, 2).Value = $row.B
    $sheetTotal.Cells.Item($r, 3).Value = $row.C
    $sheetTotal.Cells.Item($r, 4).Value = $row.D
    $r = $r + 1
}

Write-Output "done"
